$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows (bottom-up so row numbers stay valid):
# Row 12: 001759765 Natal 2188.63
# Row 10: 004500804 Rafael 5002.08
# Row 9:  004508516 Eduardo 5019.65
# Row 8:  004508504 Fernando 5023.94
# Row 7:  004886366 Renato 6565.89
$ws.Rows.Item(12).Delete()
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(7).Delete()
